$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the English-translated copy of the Element/Combination table below the original
# (rows 33-35 left blank), entering data column-by-column (A, then B, then C, then D)
# to mirror how the table was authored (matches shared-string insertion order).
# Column A
$ws.Range("A36").Value = "ElementName"
$ws.Range("A37").Value = "River"
$ws.Range("A38").Value = "Wind"
$ws.Range("A39").Value = "Earthquake"
$ws.Range("A40").Value = "Life"
$ws.Range("A41").Value = "Storm"
$ws.Range("A42").Value = "Might"
$ws.Range("A43").Value = "Mud"
$ws.Range("A44").Value = "Pressure"
$ws.Range("A45").Value = "Steam"
$ws.Range("A46").Value = "Sun"
$ws.Range("A47").Value = "Stench"
$ws.Range("A48").Value = "Fish Stick"
$ws.Range("A49").Value = "Weapon"
$ws.Range("A50").Value = "Fish"
$ws.Range("A51").Value = "Moscito"
$ws.Range("A52").Value = "Death"
$ws.Range("A53").Value = "Human"
$ws.Range("A54").Value = "Plant"
$ws.Range("A55").Value = "King"
$ws.Range("A56").Value = "War"
$ws.Range("A57").Value = "Humans"
$ws.Range("A58").Value = "Wheat"
$ws.Range("A59").Value = "Tree"
$ws.Range("A60").Value = "Zombie"
$ws.Range("A61").Value = "Fire"
$ws.Range("A62").Value = "Water"
$ws.Range("A63").Value = "Air"
$ws.Range("A64").Value = "Earth"
$ws.Range("A65").Value = "Energy"
$ws.Range("A66").Value = "Intelligence"

# Column B
$ws.Range("B36").Value = "Description"
$ws.Range("B37").Value = "A river is a natural flowing water resource."
$ws.Range("B38").Value = "Wind is air in a hurry."
$ws.Range("B39").Value = "Destructive eruption of earth."
$ws.Range("B40").Value = "Life is 42."
$ws.Range("B41").Value = "The storm is a invention of ancient meteorologists."
$ws.Range("B42").Value = "Might is the power to change the thinking of other people."
$ws.Range("B43").Value = "Mud contains 99% of commercially availble dirt."
$ws.Range("B44").Value = "Pressure is the force applied perpendicular to the surface of an object."
$ws.Range("B45").Value = "Steam consists of small water drops in the air."
$ws.Range("B46").Value = "The sun is a star which is orbiting the earth."
$ws.Range("B47").Value = "It smells… Take a deep breath."
$ws.Range("B48").Value = "A baked or fried snack similar to french fries but made of fish."
$ws.Range("B49").Value = "Als Waffe werden in der Regel alle Gegenstände bezeichnet, die Fähig sind, Lebewesen Schaden zuzufügen."
$ws.Range("B50").Value = "Fische sind aquatisch lebende Wirbeltiere mit Kiemen. "
$ws.Range("B51").Value = "Kleine, fliegende, lästige Blutsauger."
$ws.Range("B52").Value = "Der Tod ist der Zustand des Nicht-Lebens."
$ws.Range("B53").Value = "Du… Hoffentlich?!"
$ws.Range("B54").Value = "Eine Pflanze ist ein grünes Lebewesen mit Würde und Rechten."
$ws.Range("B55").Value = "Ein König ist eigentlich ein ganz normaler Mensch, nur hat er meistens eine Krone auf dem Kopf."
$ws.Range("B56").Value = "Der friedensähnliche Zustand, in welchem man sich gegenseitig umbringt und dafür gelobt wird."
$ws.Range("B57").Value = "Du… und andere, die so sind wie du sind…"
$ws.Range("B58").Value = "Getreide gehört zu den Pflanzen die auf dem Boden wachsen."
$ws.Range("B59").Value = "Bäume bestehen aus Papier und wachsen im Baumarkt."
$ws.Range("B60").Value = "Zombies sind Menschen, die jegliche Gehirnfunktionen verloren haben und vom Tod auferstanden sind."
$ws.Range("B61").Value = "Das Feuer bezeichnet die Flammenbildung bei der Verbrennung und ist ein Grundelement."
$ws.Range("B62").Value = "Wasser ist eine chemische Verbindung aus Wasserstoff und Sauerstoff und ist ein Grundelement."
$ws.Range("B63").Value = "Als Luft bezeichnet man das Gasgemisch der Erdatmosphäre. Es ist ein Grundelement."
$ws.Range("B64").Value = "Erde ist die tote organische Substanz des Bodens. Sie ist ein Grundelement."
$ws.Range("B65").Value = "Das Energievorkommen kann weder verkleinert noch vergrössert werden. Die Energie ist ein Grundelement."
$ws.Range("B66").Value = "Intelligenz bezeichnet die kognitive Leistungsfähigkeit eines Menschen. Sie ist ein Grundelement."

# Column C
$ws.Range("C36").Value = "Kombination 1"
$ws.Range("C37").Value = "Energie"
$ws.Range("C38").Value = "Energie"
$ws.Range("C39").Value = "Energie"
$ws.Range("C40").Value = "Energie"
$ws.Range("C41").Value = "Energie"
$ws.Range("C42").Value = "Energie"
$ws.Range("C43").Value = "Erde"
$ws.Range("C44").Value = "Erde"
$ws.Range("C45").Value = "Feuer"
$ws.Range("C46").Value = "Feuer"
$ws.Range("C47").Value = "Fisch"
$ws.Range("C48").Value = "Fisch"
$ws.Range("C49").Value = "Intelligenz"
$ws.Range("C50").Value = "Leben"
$ws.Range("C51").Value = "Leben"
$ws.Range("C52").Value = "Leben"
$ws.Range("C53").Value = "Leben"
$ws.Range("C54").Value = "Leben"
$ws.Range("C55").Value = "Macht"
$ws.Range("C56").Value = "Macht"
$ws.Range("C57").Value = "Mensch"
$ws.Range("C58").Value = "Pflanze"
$ws.Range("C59").Value = "Pflanze"
$ws.Range("C60").Value = "Tod"

# Column D
$ws.Range("D36").Value = "Kombination 2"
$ws.Range("D37").Value = "Wasser"
$ws.Range("D38").Value = "Luft"
$ws.Range("D39").Value = "Erde"
$ws.Range("D40").Value = "Matsch"
$ws.Range("D41").Value = "Wind"
$ws.Range("D42").Value = "Intelligenz"
$ws.Range("D43").Value = "Wasser"
$ws.Range("D44").Value = "Luft"
$ws.Range("D45").Value = "Wasser"
$ws.Range("D46").Value = "Energie"
$ws.Range("D47").Value = "Tod"
$ws.Range("D48").Value = "Feuer"
$ws.Range("D49").Value = "Feuer"
$ws.Range("D50").Value = "Wasser"
$ws.Range("D51").Value = "Luft"
$ws.Range("D52").Value = "Feuer"
$ws.Range("D53").Value = "Intelligenz"
$ws.Range("D54").Value = "Erde"
$ws.Range("D55").Value = "Mensch"
$ws.Range("D56").Value = "Waffe"
$ws.Range("D57").Value = "Mensch"
$ws.Range("D58").Value = "Sonne"
$ws.Range("D59").Value = "Energie"
$ws.Range("D60").Value = "Leben"

# Leave the cursor on the last-entered cell, as the author did when saving
[void]$ws.Range("B48").Select()